# "Season 14, matchdays prepares"
#
# The roster row for "Оксанич Кирилл" (row 15) is removed entirely from
# Sheet1. All rows below it shift up by one (rows 16-19 -> 15-18), the
# dimension shrinks from A1:W19 to A1:W18, and the now-unused shared
# string for that player's name drops out of the shared-strings table
# automatically when the workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the row first (mirrors how this was done interactively in Excel -
# clicking the row header before deleting it), then delete it, shifting
# the remaining rows up.
$row15 = $ws.Rows.Item(15)
$row15.Select() | Out-Null
$row15.Delete()
